# dictionnaire.xlsx - "category" sheet / Table2 update
#
# - Add a new mapping row: unkown -> Autres
# - Rename the table's header columns:
#     olist_category   -> main_category
#     master_category  -> master_category_fr
# - Grow Table2 (and the sheet's used range) from A1:B72 to A1:B73
# - Leave the selection on B8

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("category")
$lo = $ws.ListObjects.Item(1)

# Write the new data row first so the shared-string table picks up "unkown"
# before the renamed headers (matches the recorded edit order).
$ws.Range("A73").Value = "unkown"
$ws.Range("B73").Value = "Autres"

# Rename the header cells (this is what actually renames the ListObject's
# columns) - column B first, then column A.
$ws.Range("B1").Value = "master_category_fr"
$ws.Range("A1").Value = "main_category"

# Grow the table to cover the newly added row.
$lo.Resize($ws.Range("A1:B73"))

# Match the recorded active selection.
$ws.Range("B8").Select()
